# Introduce natural gas boiler case: add a new "Activity" block
# (heat supply, Hotel Moeschberg, 2021, natural gas boiler) right before the
# existing "energy demand, operational, Hotel Moeschberg" block, mirroring the
# structure of the existing "heat supply, Hotel Moeschberg, 2021" block
# (rows 42-52) which uses a wood pellet furnace instead of a gas boiler.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room: insert 11 new (blank) rows above the old row 54, pushing the
#    "energy demand, operational, Hotel Moeschberg" block (old rows 54-64)
#    down to new rows 65-75. Excel automatically rewrites the relative
#    formulas that pointed into/around the shifted block.
$ws.Rows("54:64").Insert()

# 2. Fill in the new "heat supply, Hotel Moeschberg, 2021, natural gas boiler"
#    activity block in rows 54-63 (row 64 is left blank, as a separator,
#    exactly like row 53 before it and row 64 after the relocated block).

# Row 54: Activity header (bold, 12pt - same style as other "Activity" rows)
$ws.Range("A54").Value = "Activity"
$ws.Range("B54").Value = "heat supply, Hotel Moeschberg, 2021, natural gas boiler"

# Row 55: location
$ws.Range("A55").Value = "location"
$ws.Range("B55").Value = "CH"

# Row 56: reference product
$ws.Range("A56").Value = "reference product"
$ws.Range("B56").Value = "heat, Hotel Moeschberg"

# Row 57: type
$ws.Range("A57").Value = "type"
$ws.Range("B57").Value = "process"

# Row 58: unit
$ws.Range("A58").Value = "unit"
$ws.Range("B58").Value = "megajoule"

# Row 59: Exchanges header (bold, 12pt)
$ws.Range("A59").Value = "Exchanges"

# Row 60: column headers for the exchange table
$ws.Range("A60").Value = "name"
$ws.Range("B60").Value = "amount"
$ws.Range("C60").Value = "database"
$ws.Range("D60").Value = "location"
$ws.Range("E60").Value = "unit"
$ws.Range("F60").Value = "type"
$ws.Range("G60").Value = "reference product"

# Row 61: the activity's own production exchange, referencing the fields
# defined above via formulas (same pattern as row 50 for the other heat
# supply activity).
$ws.Range("A61").Formula = "=B54"
$ws.Range("B61").Value = 1
$ws.Range("C61").Formula = "=`$B`$1"
$ws.Range("D61").Formula = "=B55"
$ws.Range("E61").Formula = "=B58"
$ws.Range("F61").Value = "production"
$ws.Range("G61").Formula = "=B56"

# Row 62: existing solar collector technosphere exchange (same values as the
# other heat supply activity's matching row).
$ws.Range("A62").Value = "operation, solar collector system, Cu flat plate collector, one-family house, for hot water"
$ws.Range("B62").Value = 0.004668
$ws.Range("C62").Value = "ei 3.8 cutoff"
$ws.Range("D62").Value = "CH"
$ws.Range("E62").Value = "megajoule"
$ws.Range("F62").Value = "technosphere"
$ws.Range("G62").Value = "heat, central or small-scale, other than natural gas"

# Row 63: new natural gas boiler technosphere exchange (uses a smaller, 10pt
# font - a style not used anywhere else in the workbook yet).
$ws.Range("A63").Value = "heat production, natural gas, at boiler condensing modulating <100kW"
$ws.Range("B63").Value = 0.9953
$ws.Range("C63").Value = "ei 3.8 cutoff"
$ws.Range("D63").Value = "CH"
$ws.Range("E63").Value = "megajoule"
$ws.Range("F63").Value = "technosphere"
$ws.Range("G63").Value = "heat, central or small-scale, natural gas"

# 3. Formatting.
#    a) Copy the existing bold/12pt "section header" formatting (used by the
#       other "Activity"/"Exchanges" headers) onto the new header cells,
#       reusing the existing style instead of creating a new one.
$ws.Range("A42:B42").Copy() | Out-Null
$ws.Range("A54:B54").PasteSpecial(-4122) | Out-Null
$ws.Range("A48").Copy() | Out-Null
$ws.Range("A59").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

#    b) Give rows 54 and 59 the taller row height used by other section
#       header rows.
$ws.Rows(54).RowHeight = 15.6
$ws.Rows(59).RowHeight = 15.6

#    c) Create (once) a 10pt font style on a scratch cell and copy that
#       formatting onto A63/G63, so only a single new font/style entry is
#       added to the workbook instead of one per cell.
$scratch = $ws.Range("ZZ1")
$scratch.Value = "x"
$scratch.Font.Size = 10
$scratch.HorizontalAlignment = 1
$scratch.Copy() | Out-Null
$ws.Range("A63").PasteSpecial(-4122) | Out-Null
$ws.Range("G63").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$scratch.Clear() | Out-Null

# 4. Update the view state to match: scroll so row 43 is at the top and
#    select D58.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1
$ws.Range("D58").Select()
